$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the Salary VR formula (typo 1.1254 -> 1.254)
$ws.Range("B2").Formula = "=200000*1.254*3"

# New "Comments" column (C) with explanatory notes
$ws.Range("C1").Value = "Comments"
$ws.Range("C2").Value = "15% of salary per year for 3 years"
$ws.Range("C3").Value = "Project grant from SLS"
$ws.Range("C4").Value = "Funding available for project costs from ALF grant "
$ws.Range("C5").Value = "Remaining funds from 7550000 VR grant after salary and planned invoices from TGI for batch 1 and 2 have been accounted for"

# New rows 9-10 (row 8 left blank)
$ws.Range("A9").Value = "Total salary remaining"
$ws.Range("B9").Formula = "=510000*5-B2-(B2/3*2)"

$ws.Range("A10").Value = "Total available if salary included "
$ws.Range("B10").Formula = "=(B6+B9)/1.254"

$ws.Range("F9").Select()
